# Replace the "black" (⬛ / noir) status entries with "blue" (📘 / bleu)
# across the whole data range of the sheet.
#
# Column A = statut, Column B = statut_label.
# Every data row (2..20) currently has A="⬛" and B="noir"; these are
# stored as shared strings, so editing their text effectively updates
# every cell that references them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    if ($cellA.Value2 -eq "⬛") {
        $cellA.Value = "📘"
    }
    if ($cellB.Value2 -eq "noir") {
        $cellB.Value = "bleu"
    }
}
